# Updates Leve-profit market-price figures (columns H-N) across all job sheets.
# Values correspond to a scheduled market-data refresh; see commit message.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 296318.8
$ws.Range("I6").Value = 370175.5
$ws.Range("J6").Value = 892
$ws.Range("K6").Value = 1110526.5
$ws.Range("L6").Value = 2676
$ws.Range("M6").Value = -1110414.5
$ws.Range("N6").Value = -2900
# Row 55
$ws.Range("H55").Value = 101.333336
$ws.Range("J55").Value = 116.5
$ws.Range("L55").Value = 116.5
$ws.Range("N55").Value = -544.5
# Row 87
$ws.Range("H87").Value = 30522.857
$ws.Range("J87").Value = 30522.857
$ws.Range("L87").Value = 30522.857
$ws.Range("N87").Value = -33018.857
# Row 90
$ws.Range("H90").Value = 30522.857
$ws.Range("J90").Value = 30522.857
$ws.Range("L90").Value = 91568.571
$ws.Range("N90").Value = -104048.571
# Row 136
$ws.Range("H136").Value = 53296.668
$ws.Range("J136").Value = 53296.668
$ws.Range("L136").Value = 53296.668
$ws.Range("N136").Value = -63496.668
# Row 137
$ws.Range("H137").Value = 1905.1818
$ws.Range("I137").Value = 2077.6365
$ws.Range("K137").Value = 6232.9095
$ws.Range("M137").Value = -3682.9095

$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 2854.625
$ws.Range("I26").Value = 2088.5454
$ws.Range("J26").Value = 4540
$ws.Range("K26").Value = 2088.5454
$ws.Range("L26").Value = 4540
$ws.Range("M26").Value = -1758.5454
$ws.Range("N26").Value = -5200
# Row 32
$ws.Range("H32").Value = 27482.338
$ws.Range("I32").Value = 6235.3066
$ws.Range("K32").Value = 6235.3066
$ws.Range("M32").Value = -5948.3066
# Row 44
$ws.Range("H44").Value = 11661.625
$ws.Range("J44").Value = 12892.714
$ws.Range("L44").Value = 12892.714
$ws.Range("N44").Value = -13868.714
# Row 55
$ws.Range("H55").Value = 12237
$ws.Range("J55").Value = 12891.625
$ws.Range("L55").Value = 12891.625
$ws.Range("N55").Value = -13521.625
# Row 61
$ws.Range("H61").Value = 2424.9167
$ws.Range("I61").Value = 2319.9
$ws.Range("K61").Value = 2319.9
$ws.Range("M61").Value = -2107.9
# Row 80
$ws.Range("H80").Value = 27157.666
$ws.Range("J80").Value = 27157.666
$ws.Range("L80").Value = 27157.666
$ws.Range("N80").Value = -29153.666
# Row 83
$ws.Range("H83").Value = 27157.666
$ws.Range("J83").Value = 27157.666
$ws.Range("L83").Value = 81472.99800000001
$ws.Range("N83").Value = -91456.99800000001
# Row 122
$ws.Range("H122").Value = 1752.2778
$ws.Range("I122").Value = 2026.1
$ws.Range("K122").Value = 6078.299999999999
$ws.Range("M122").Value = -3628.299999999999
# Row 132
$ws.Range("H132").Value = 4737.8335
$ws.Range("I132").Value = 5318
$ws.Range("K132").Value = 15954
$ws.Range("M132").Value = -13424
# Row 136
$ws.Range("H136").Value = 2424.9167
$ws.Range("I136").Value = 2319.9
$ws.Range("K136").Value = 6959.700000000001
$ws.Range("M136").Value = -4409.700000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 424.38235
$ws.Range("I94").Value = 387.85715
$ws.Range("J94").Value = 483.3846
$ws.Range("K94").Value = 387.85715
$ws.Range("L94").Value = 483.3846
$ws.Range("M94").Value = 63.14285000000001
$ws.Range("N94").Value = -1385.3846
# Row 105
$ws.Range("H105").Value = 252415
$ws.Range("I105").Value = 169518.17
$ws.Range("K105").Value = 169518.17
$ws.Range("M105").Value = -167771.17

$ws = $wb.Worksheets.Item("CRP")
# Row 35
$ws.Range("H35").Value = 11500
$ws.Range("I35").Value = 11000
$ws.Range("J35").Value = 12000
$ws.Range("K35").Value = 11000
$ws.Range("L35").Value = 12000
$ws.Range("M35").Value = -10706
$ws.Range("N35").Value = -12588
# Row 134
$ws.Range("H134").Value = 2333.3333
$ws.Range("I134").Value = 1500
$ws.Range("K134").Value = 4500
$ws.Range("M134").Value = -1965

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 250.5
$ws.Range("J7").Value = 250.5
$ws.Range("L7").Value = 751.5
$ws.Range("N7").Value = -975.5
# Row 80
$ws.Range("H80").Value = 11040.3
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 12155.889
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 36467.667
$ws.Range("M80").Value = -2064
$ws.Range("N80").Value = -38339.667
# Row 83
$ws.Range("H83").Value = 11040.3
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 12155.889
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 109403.001
$ws.Range("M83").Value = -4320
$ws.Range("N83").Value = -118763.001
# Row 92
$ws.Range("H92").Value = 659.8
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 899.5
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 2698.5
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -5194.5
# Row 133
$ws.Range("H133").Value = 3200
$ws.Range("I133").Value = 3200
$ws.Range("K133").Value = 9600
$ws.Range("M133").Value = -4540

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1883.1666
$ws.Range("J113").Value = 1883.1666
$ws.Range("L113").Value = 1883.1666
$ws.Range("N113").Value = -6223.1666
# Row 126
$ws.Range("H126").Value = 2570.5715
$ws.Range("I126").Value = 3045
$ws.Range("J126").Value = 1938
$ws.Range("K126").Value = 9135
$ws.Range("L126").Value = 5814
$ws.Range("M126").Value = -6665
$ws.Range("N126").Value = -10754

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3794
$ws.Range("I68").Value = 3484.5
$ws.Range("J68").Value = 3979.7
$ws.Range("K68").Value = 3484.5
$ws.Range("L68").Value = 3979.7
$ws.Range("M68").Value = -2735.5
$ws.Range("N68").Value = -5477.7
# Row 71
$ws.Range("H71").Value = 3794
$ws.Range("I71").Value = 3484.5
$ws.Range("J71").Value = 3979.7
$ws.Range("K71").Value = 17422.5
$ws.Range("L71").Value = 19898.5
$ws.Range("M71").Value = -13678.5
$ws.Range("N71").Value = -27386.5

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4810654.5
$ws.Range("I62").Value = 25645358
$ws.Range("J62").Value = 2646
$ws.Range("K62").Value = 25645358
$ws.Range("L62").Value = 2646
$ws.Range("M62").Value = -25644734
$ws.Range("N62").Value = -3894
# Row 65
$ws.Range("H65").Value = 4810654.5
$ws.Range("I65").Value = 25645358
$ws.Range("J65").Value = 2646
$ws.Range("K65").Value = 128226790
$ws.Range("L65").Value = 13230
$ws.Range("M65").Value = -128223670
$ws.Range("N65").Value = -19470
# Row 81
$ws.Range("H81").Value = 137883.25
$ws.Range("J81").Value = 549999.5
$ws.Range("L81").Value = 1099999
$ws.Range("N81").Value = -1102121
# Row 84
$ws.Range("H84").Value = 137883.25
$ws.Range("J84").Value = 549999.5
$ws.Range("L84").Value = 5499995
$ws.Range("N84").Value = -5510603
# Row 126
$ws.Range("H126").Value = 1331.5
$ws.Range("I126").Value = 1205
$ws.Range("J126").Value = 1559.2
$ws.Range("K126").Value = 3615
$ws.Range("L126").Value = 4677.6
$ws.Range("M126").Value = -1145
$ws.Range("N126").Value = -9617.6

